# Update "想去人数" (want-to-go count) figures in column F across sheets.
# Source data refresh (gh-pages regeneration) bumped several event counts.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsLocal   = $wb.Worksheets.Item("本地生活")
$wsAll     = $wb.Worksheets.Item("全部类型")

# 展览 ("Exhibition") sheet
$wsExhibit.Range("F10").Value = 590
$wsExhibit.Range("F12").Value = 339
$wsExhibit.Range("F14").Value = 6029
$wsExhibit.Range("F15").Value = 646
$wsExhibit.Range("F16").Value = 1055
$wsExhibit.Range("F21").Value = 570
$wsExhibit.Range("F23").Value = 50
$wsExhibit.Range("F25").Value = 138
$wsExhibit.Range("F28").Value = 1009
$wsExhibit.Range("F32").Value = 14
$wsExhibit.Range("F35").Value = 3337

# 本地生活 ("Local life") sheet
$wsLocal.Range("F6").Value = 1149

# 全部类型 ("All types") sheet
$wsAll.Range("F5").Value = 1149
$wsAll.Range("F14").Value = 590
$wsAll.Range("F17").Value = 339
$wsAll.Range("F19").Value = 6029
$wsAll.Range("F21").Value = 646
$wsAll.Range("F22").Value = 1055
$wsAll.Range("F27").Value = 570
$wsAll.Range("F39").Value = 1009
$wsAll.Range("F49").Value = 3337
